# Update TPM-derived NATMI LR-pair statistics for Inhba-Acvr1b
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 0.08097566666666667
$ws.Range("H2").Value = 0.242927
$ws.Range("I2").Value = 0.005588990034505014
$ws.Range("J2").Value = 0.005588990034505015
$ws.Range("M2").Value = 1.923239
$ws.Range("N2").Value = 5.769717
$ws.Range("O2").Value = 0.2340262838603868
$ws.Range("P2").Value = 0.2340262838603868
$ws.Range("Q2").Value = 0.1557355601843333
$ws.Range("R2").Value = 1.401620041659
$ws.Range("S2").Value = 0.001307970568307944
$ws.Range("T2").Value = 0.001307970568307944

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 0.08097566666666667
$ws.Range("H3").Value = 0.242927
$ws.Range("I3").Value = 0.005588990034505014
$ws.Range("J3").Value = 0.005588990034505015
$ws.Range("O3").Value = 0.4335574295612247
$ws.Range("P3").Value = 0.4335574295612246
$ws.Range("Q3").Value = 0.2885159224468889
$ws.Range("R3").Value = 2.596643302022
$ws.Range("S3").Value = 0.002423148153203294
$ws.Range("T3").Value = 0.002423148153203294

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 0.08097566666666667
$ws.Range("H4").Value = 0.242927
$ws.Range("I4").Value = 0.005588990034505014
$ws.Range("J4").Value = 0.005588990034505015
$ws.Range("O4").Value = 0.3324162865783886
$ws.Range("P4").Value = 0.3324162865783886
$ws.Range("Q4").Value = 0.2212103518917778
$ws.Range("R4").Value = 1.990893167026
$ws.Range("S4").Value = 0.001857871312993777
$ws.Range("T4").Value = 0.001857871312993777

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.6976944377922635
$ws.Range("J5").Value = 0.6976944377922635
$ws.Range("M5").Value = 1.923239
$ws.Range("N5").Value = 5.769717
$ws.Range("O5").Value = 0.2340262838603868
$ws.Range("P5").Value = 0.2340262838603868
$ws.Range("Q5").Value = 19.441049892066
$ws.Range("R5").Value = 174.969449028594
$ws.Range("S5").Value = 0.1632788365465852
$ws.Range("T5").Value = 0.1632788365465852

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.6976944377922635
$ws.Range("J6").Value = 0.6976944377922635
$ws.Range("O6").Value = 0.4335574295612247
$ws.Range("P6").Value = 0.4335574295612246
$ws.Range("S6").Value = 0.3024906070683775
$ws.Range("T6").Value = 0.3024906070683775

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.6976944377922635
$ws.Range("J7").Value = 0.6976944377922635
$ws.Range("O7").Value = 0.3324162865783886
$ws.Range("P7").Value = 0.3324162865783886
$ws.Range("S7").Value = 0.2319249941773008
$ws.Range("T7").Value = 0.2319249941773008

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.2967165721732315
$ws.Range("J8").Value = 0.2967165721732316
$ws.Range("M8").Value = 1.923239
$ws.Range("N8").Value = 5.769717
$ws.Range("O8").Value = 0.2340262838603868
$ws.Range("P8").Value = 0.2340262838603868
$ws.Range("Q8").Value = 8.267919838484001
$ws.Range("R8").Value = 74.411278546356
$ws.Range("S8").Value = 0.06943947674549363
$ws.Range("T8").Value = 0.06943947674549363

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.2967165721732315
$ws.Range("J9").Value = 0.2967165721732316
$ws.Range("O9").Value = 0.4335574295612247
$ws.Range("P9").Value = 0.4335574295612246
$ws.Range("S9").Value = 0.1286436743396439
$ws.Range("T9").Value = 0.1286436743396439

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.2967165721732315
$ws.Range("J10").Value = 0.2967165721732316
$ws.Range("O10").Value = 0.3324162865783886
$ws.Range("P10").Value = 0.3324162865783886
$ws.Range("S10").Value = 0.09863342108809404
$ws.Range("T10").Value = 0.09863342108809406
